## Categories sheet: the "quarterly variables" C column held a hard-coded
## boolean FALSE for every data row (C2:C301). The fix re-enters each of
## those cells as the literal text "FALSE" instead of the Boolean FALSE,
## so the column becomes a text/category column like the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Categories")

# Select column C first (matches how this was done interactively - the
# whole column was selected before retyping the values).
$ws.Columns.Item(3).Select() | Out-Null

$rng = $ws.Range("C2:C301")

# Writing the literal string "FALSE" straight into .Value/.Value2 gets
# auto-coerced back to a Boolean by the COM layer (same as typing FALSE
# into Excel without a leading apostrophe). Instead, push it in through a
# text formula (forces a String result) and then paste-special just the
# values back on top of themselves - this "bakes" the formula result into
# a plain text cell (t="s") without leaving a formula behind and without
# picking up a quote-prefix / "number stored as text" style variant.
$rng.Formula = "=""FALSE"""
$rng.Copy() | Out-Null
$rng.PasteSpecial(-4163) | Out-Null
